$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with new content
$ws.Range("A2").Value = "Ecalation Access"
$ws.Range("B2").Value = "User need Escalation approval where Approve/Reject button is visible in case approval history"
$ws.Range("C2").Value = "We need to add the user in the Case Escalation Request Queue"
$ws.Range("D2").Value = "Request access"

# Overwrite row 3 with new content
$ws.Range("A3").Value = "Credit Hold"
$ws.Range("B3").Value = "Customer not able to create cases from SC"
$ws.Range("C3").Value = "Send and email to DLcollections@csod.com."
$ws.Range("D3").Value = "Question"

# Add new row 4
$ws.Range("A4").Value = "Not able to add DSC"
$ws.Range("B4").Value = 'Getting this error "This contact has the same email address as an active Support Central user. Duplicate Support Central users cannot be created, so this contact cannot be added as a Designated Support Contact. Please contact IT Support for assistance."'
$ws.Range("C4").Value = "Search the contact from User in salesforce and add .Invalid/Test at the end of the email and save"
$ws.Range("D4").Value = "Issue"

# Set column widths to reflect new content (matching Excel's auto-fit results)
$ws.Columns.Item(1).ColumnWidth = 12.833333333333334
$ws.Columns.Item(2).ColumnWidth = 17.333333333333332
$ws.Columns.Item(3).ColumnWidth = 52.333333333333336
$ws.Columns.Item(4).ColumnWidth = 4.5

# Update selection to match target (A5)
$ws.Range("A5").Select() | Out-Null
